$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply date format (style) to column B for new rows by copying format from B283,
# which keeps the existing date-format style index instead of creating a new one.
$ws.Range("B283").Copy() | Out-Null
$ws.Range("B284:B299").PasteSpecial(-4122) | Out-Null

# Row 284 (G:V) uses centered horizontal+vertical alignment (matches existing style index 2).
$ws.Range("G284:V284").HorizontalAlignment = -4108
$ws.Range("G284:V284").VerticalAlignment = -4108

# E296 and E297 reuse the existing "alt font" style (style index 6) already used elsewhere
# in column E (e.g. E245:E248, E278:E280); copy format from E245.
$ws.Range("E245").Copy() | Out-Null
$ws.Range("E296").PasteSpecial(-4122) | Out-Null
$ws.Range("E245").Copy() | Out-Null
$ws.Range("E297").PasteSpecial(-4122) | Out-Null

# Row 284
$ws.Range("A284").Value = "Entrainement"
$ws.Range("B284").Value = 45880
$ws.Range("C284").Value = "Global"
$ws.Range("E284").Value = "Hedi Nasri"
$ws.Range("F284").Value = "right back"
$ws.Range("G284").Value = "00:52:10"
$ws.Range("H284").Value = 4.62
$ws.Range("I284").Value = 0.1
$ws.Range("J284").Value = 4.5199999999999996
$ws.Range("K284").Value = 0.09
$ws.Range("L284").Value = 0.02
$ws.Range("M284").Value = 0
$ws.Range("N284").Value = 0
$ws.Range("O284").Value = 0
$ws.Range("P284").Value = 5.26
$ws.Range("Q284").Value = 24.79
$ws.Range("R284").Value = 4.8099999999999996
$ws.Range("S284").Value = 12
$ws.Range("T284").Value = 1
$ws.Range("U284").Value = 13
$ws.Range("V284").Value = 3

# Row 285
$ws.Range("A285").Value = "Entrainement"
$ws.Range("B285").Value = 45881
$ws.Range("C285").Value = "Global"
$ws.Range("E285").Value = "Maé Clavel"
$ws.Range("F285").Value = "left back"
$ws.Range("G285").Value = "00:48:29"
$ws.Range("H285").Value = 5.83
$ws.Range("I285").Value = 1.2
$ws.Range("J285").Value = 4.6100000000000003
$ws.Range("K285").Value = 0.84
$ws.Range("L285").Value = 0.33
$ws.Range("M285").Value = 0.04
$ws.Range("N285").Value = 0
$ws.Range("O285").Value = 4
$ws.Range("P285").Value = 7.22
$ws.Range("Q285").Value = 27.04
$ws.Range("R285").Value = 4.51
$ws.Range("S285").Value = 20
$ws.Range("T285").Value = 3
$ws.Range("U285").Value = 11
$ws.Range("V285").Value = 10

# Row 286
$ws.Range("A286").Value = "Entrainement"
$ws.Range("B286").Value = 45881
$ws.Range("C286").Value = "Global"
$ws.Range("E286").Value = "Levy Ndoutoume"
$ws.Range("F286").Value = "left back"
$ws.Range("G286").Value = "00:50:12"
$ws.Range("H286").Value = 5.14
$ws.Range("I286").Value = 0.85
$ws.Range("J286").Value = 4.28
$ws.Range("K286").Value = 0.57999999999999996
$ws.Range("L286").Value = 0.21
$ws.Range("M286").Value = 0.08
$ws.Range("N286").Value = 0
$ws.Range("O286").Value = 4
$ws.Range("P286").Value = 6.13
$ws.Range("Q286").Value = 29.96
$ws.Range("R286").Value = 4.6100000000000003
$ws.Range("S286").Value = 20
$ws.Range("T286").Value = 5
$ws.Range("U286").Value = 23
$ws.Range("V286").Value = 9

# Row 287
$ws.Range("A287").Value = "Entrainement"
$ws.Range("B287").Value = 45881
$ws.Range("C287").Value = "Global"
$ws.Range("E287").Value = "Ilan Ihaddadene"
$ws.Range("F287").Value = "center midfield"
$ws.Range("G287").Value = "00:48:29"
$ws.Range("H287").Value = 5.75
$ws.Range("I287").Value = 1.04
$ws.Range("J287").Value = 4.7
$ws.Range("K287").Value = 0.9
$ws.Range("L287").Value = 0.15
$ws.Range("M287").Value = 0
$ws.Range("N287").Value = 0
$ws.Range("O287").Value = 0
$ws.Range("P287").Value = 7.1
$ws.Range("Q287").Value = 24.73
$ws.Range("R287").Value = 4.51
$ws.Range("S287").Value = 21
$ws.Range("T287").Value = 7
$ws.Range("U287").Value = 21
$ws.Range("V287").Value = 8

# Row 288
$ws.Range("A288").Value = "Entrainement"
$ws.Range("B288").Value = 45881
$ws.Range("C288").Value = "Global"
$ws.Range("E288").Value = "Amir Etien"
$ws.Range("F288").Value = "right forward"
$ws.Range("G288").Value = "01:12:07"
$ws.Range("H288").Value = 6.81
$ws.Range("I288").Value = 1.26
$ws.Range("J288").Value = 5.53
$ws.Range("K288").Value = 0.78
$ws.Range("L288").Value = 0.36
$ws.Range("M288").Value = 0.14000000000000001
$ws.Range("N288").Value = 0
$ws.Range("O288").Value = 9
$ws.Range("P288").Value = 5.64
$ws.Range("Q288").Value = 29.8
$ws.Range("R288").Value = 4.4400000000000004
$ws.Range("S288").Value = 35
$ws.Range("T288").Value = 7
$ws.Range("U288").Value = 24
$ws.Range("V288").Value = 10

# Row 289
$ws.Range("A289").Value = "Entrainement"
$ws.Range("B289").Value = 45881
$ws.Range("C289").Value = "Global"
$ws.Range("E289").Value = "Naim Dhib"
$ws.Range("F289").Value = "center midfield"
$ws.Range("G289").Value = "01:20:13"
$ws.Range("H289").Value = 8.25
$ws.Range("I289").Value = 1.58
$ws.Range("J289").Value = 6.65
$ws.Range("K289").Value = 1.1499999999999999
$ws.Range("L289").Value = 0.38
$ws.Range("M289").Value = 0.07
$ws.Range("N289").Value = 0
$ws.Range("O289").Value = 7
$ws.Range("P289").Value = 6.14
$ws.Range("Q289").Value = 29.94
$ws.Range("R289").Value = 4.75
$ws.Range("S289").Value = 37
$ws.Range("T289").Value = 6
$ws.Range("U289").Value = 32
$ws.Range("V289").Value = 13

# Row 290
$ws.Range("A290").Value = "Entrainement"
$ws.Range("B290").Value = 45881
$ws.Range("C290").Value = "Global"
$ws.Range("E290").Value = "Yanis Berrached"
$ws.Range("F290").Value = "center midfield"
$ws.Range("G290").Value = "01:14:46"
$ws.Range("H290").Value = 8.09
$ws.Range("I290").Value = 1.37
$ws.Range("J290").Value = 6.71
$ws.Range("K290").Value = 1.06
$ws.Range("L290").Value = 0.26
$ws.Range("M290").Value = 0.06
$ws.Range("N290").Value = 0
$ws.Range("O290").Value = 6
$ws.Range("P290").Value = 6.55
$ws.Range("Q290").Value = 29.09
$ws.Range("R290").Value = 4.37
$ws.Range("S290").Value = 20
$ws.Range("T290").Value = 4
$ws.Range("U290").Value = 17
$ws.Range("V290").Value = 7

# Row 291
$ws.Range("A291").Value = "Entrainement"
$ws.Range("B291").Value = 45881
$ws.Range("C291").Value = "Global"
$ws.Range("E291").Value = "Amine Taiar"
$ws.Range("F291").Value = "center back"
$ws.Range("G291").Value = "00:50:12"
$ws.Range("H291").Value = 5.55
$ws.Range("I291").Value = 0.94
$ws.Range("J291").Value = 4.5999999999999996
$ws.Range("K291").Value = 0.64
$ws.Range("L291").Value = 0.22
$ws.Range("M291").Value = 0.08
$ws.Range("N291").Value = 0
$ws.Range("O291").Value = 8
$ws.Range("P291").Value = 6.6
$ws.Range("Q291").Value = 28.9
$ws.Range("R291").Value = 4.22
$ws.Range("S291").Value = 13
$ws.Range("T291").Value = 3
$ws.Range("U291").Value = 22
$ws.Range("V291").Value = 1

# Row 292
$ws.Range("A292").Value = "Entrainement"
$ws.Range("B292").Value = 45881
$ws.Range("C292").Value = "Global"
$ws.Range("E292").Value = "Naim Ighbane"
$ws.Range("F292").Value = "center back"
$ws.Range("G292").Value = "01:12:07"
$ws.Range("H292").Value = 7.15
$ws.Range("I292").Value = 1.01
$ws.Range("J292").Value = 6.13
$ws.Range("K292").Value = 0.76
$ws.Range("L292").Value = 0.2
$ws.Range("M292").Value = 0.06
$ws.Range("N292").Value = 0
$ws.Range("O292").Value = 6
$ws.Range("P292").Value = 5.92
$ws.Range("Q292").Value = 28.67
$ws.Range("R292").Value = 4.33
$ws.Range("S292").Value = 23
$ws.Range("T292").Value = 2
$ws.Range("U292").Value = 15
$ws.Range("V292").Value = 4

# Row 293
$ws.Range("A293").Value = "Entrainement"
$ws.Range("B293").Value = 45881
$ws.Range("C293").Value = "Global"
$ws.Range("E293").Value = "Karahali Souaré"
$ws.Range("F293").Value = "right forward"
$ws.Range("G293").Value = "01:12:39"
$ws.Range("H293").Value = 8.11
$ws.Range("I293").Value = 1.56
$ws.Range("J293").Value = 6.52
$ws.Range("K293").Value = 1.04
$ws.Range("L293").Value = 0.46
$ws.Range("M293").Value = 0.08
$ws.Range("N293").Value = 0
$ws.Range("O293").Value = 8
$ws.Range("P293").Value = 6.67
$ws.Range("Q293").Value = 29.28
$ws.Range("R293").Value = 4.76
$ws.Range("S293").Value = 50
$ws.Range("T293").Value = 8
$ws.Range("U293").Value = 48
$ws.Range("V293").Value = 16

# Row 294
$ws.Range("A294").Value = "Entrainement"
$ws.Range("B294").Value = 45881
$ws.Range("C294").Value = "Global"
$ws.Range("E294").Value = "Romain Thunet"
$ws.Range("F294").Value = "center back"
$ws.Range("G294").Value = "00:48:29"
$ws.Range("H294").Value = 5.25
$ws.Range("I294").Value = 0.85
$ws.Range("J294").Value = 4.38
$ws.Range("K294").Value = 0.54
$ws.Range("L294").Value = 0.28999999999999998
$ws.Range("M294").Value = 0.04
$ws.Range("N294").Value = 0
$ws.Range("O294").Value = 7
$ws.Range("P294").Value = 6.53
$ws.Range("Q294").Value = 27.25
$ws.Range("R294").Value = 4.45
$ws.Range("S294").Value = 21
$ws.Range("T294").Value = 3
$ws.Range("U294").Value = 14
$ws.Range("V294").Value = 3

# Row 295
$ws.Range("A295").Value = "Entrainement"
$ws.Range("B295").Value = 45881
$ws.Range("C295").Value = "Global"
$ws.Range("E295").Value = "Mattheo Haon"
$ws.Range("F295").Value = "right back"
$ws.Range("G295").Value = "01:12:07"
$ws.Range("H295").Value = 8
$ws.Range("I295").Value = 1.34
$ws.Range("J295").Value = 6.64
$ws.Range("K295").Value = 0.85
$ws.Range("L295").Value = 0.38
$ws.Range("M295").Value = 0.13
$ws.Range("N295").Value = 0
$ws.Range("O295").Value = 10
$ws.Range("P295").Value = 6.63
$ws.Range("Q295").Value = 30.17
$ws.Range("R295").Value = 4.76
$ws.Range("S295").Value = 17
$ws.Range("T295").Value = 10
$ws.Range("U295").Value = 16
$ws.Range("V295").Value = 7

# Row 296
$ws.Range("A296").Value = "Entrainement"
$ws.Range("B296").Value = 45881
$ws.Range("C296").Value = "Global"
$ws.Range("E296").Value = "Kamal Bafounta"
$ws.Range("F296").Value = "left forward"
$ws.Range("G296").Value = "00:25:14"
$ws.Range("H296").Value = 2.5299999999999998
$ws.Range("I296").Value = 0.44
$ws.Range("J296").Value = 2.09
$ws.Range("K296").Value = 0.28000000000000003
$ws.Range("L296").Value = 0.12
$ws.Range("M296").Value = 0.04
$ws.Range("N296").Value = 0
$ws.Range("O296").Value = 3
$ws.Range("P296").Value = 6.01
$ws.Range("Q296").Value = 30.36
$ws.Range("R296").Value = 4.37
$ws.Range("S296").Value = 11
$ws.Range("T296").Value = 1
$ws.Range("U296").Value = 9
$ws.Range("V296").Value = 4

# Row 297
$ws.Range("A297").Value = "Entrainement"
$ws.Range("B297").Value = 45881
$ws.Range("C297").Value = "Global"
$ws.Range("E297").Value = "Omar Benyounes"
$ws.Range("F297").Value = "center midfield"
$ws.Range("G297").Value = "00:25:54"
$ws.Range("H297").Value = 2.84
$ws.Range("I297").Value = 0.67
$ws.Range("J297").Value = 2.16
$ws.Range("K297").Value = 0.47
$ws.Range("L297").Value = 0.21
$ws.Range("M297").Value = 0
$ws.Range("N297").Value = 0
$ws.Range("O297").Value = 1
$ws.Range("P297").Value = 6.55
$ws.Range("Q297").Value = 25.09
$ws.Range("R297").Value = 4.3600000000000003
$ws.Range("S297").Value = 10
$ws.Range("T297").Value = 3
$ws.Range("U297").Value = 9
$ws.Range("V297").Value = 9

# Row 298
$ws.Range("A298").Value = "Entrainement"
$ws.Range("B298").Value = 45881
$ws.Range("C298").Value = "Global"
$ws.Range("E298").Value = "Emmanuel Valey"
$ws.Range("F298").Value = "left forward"
$ws.Range("G298").Value = "00:50:12"
$ws.Range("H298").Value = 5.57
$ws.Range("I298").Value = 1
$ws.Range("J298").Value = 4.55
$ws.Range("K298").Value = 0.65
$ws.Range("L298").Value = 0.32
$ws.Range("M298").Value = 0.05
$ws.Range("N298").Value = 0
$ws.Range("O298").Value = 6
$ws.Range("P298").Value = 6.63
$ws.Range("Q298").Value = 28.6
$ws.Range("R298").Value = 4.57
$ws.Range("S298").Value = 17
$ws.Range("T298").Value = 3
$ws.Range("U298").Value = 19
$ws.Range("V298").Value = 7

# Row 299
$ws.Range("A299").Value = "Entrainement"
$ws.Range("B299").Value = 45881
$ws.Range("C299").Value = "Global"
$ws.Range("E299").Value = "Rayane Chayebi"
$ws.Range("F299").Value = "center midfield"
$ws.Range("G299").Value = "01:12:22"
$ws.Range("H299").Value = 7.41
$ws.Range("I299").Value = 1
$ws.Range("J299").Value = 6.4
$ws.Range("K299").Value = 0.84
$ws.Range("L299").Value = 0.15
$ws.Range("M299").Value = 0.03
$ws.Range("N299").Value = 0
$ws.Range("O299").Value = 4
$ws.Range("P299").Value = 6.12
$ws.Range("Q299").Value = 26.96
$ws.Range("R299").Value = 5.01
$ws.Range("S299").Value = 21
$ws.Range("T299").Value = 9
$ws.Range("U299").Value = 28
$ws.Range("V299").Value = 5

# Update the active selection to reflect the new end of the data (matches the
# author continuing to enter rows after this point).
$ws.Range("C303").Select()
